# The commit swaps the "Integral" (Red Violet) theme and the plain
# "Office Theme" between the deck's two theme parts: the presentation's
# real/applied theme goes from Red Violet -> Office colors, while the
# (otherwise unused) secondary theme part picks up the Red Violet colors.
#
# Concretely, for the applied theme (the one driving SlideMaster /
# Presentation), every theme color slot is repointed from the "Integral"
# palette to the stock "Office" palette:
#   dk1      000000 -> 000000
#   lt1      FFFFFF -> FFFFFF
#   dk2      454551 -> 44546A
#   lt2      D8D9DC -> E7E6E6
#   accent1  E32D91 -> 5B9BD5
#   accent2  C830CC -> ED7D31
#   accent3  4EA6DC -> A5A5A5
#   accent4  4775E7 -> FFC000
#   accent5  8971E1 -> 4472C4
#   accent6  D54773 -> 70AD47
#   hlink    6B9F25 -> 0563C1
#   folHlink 8C8C8C -> 954F72

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$scheme = $master.ColorScheme

# ppColorSchemeIndex slots 1-8 (dk1, lt1, dk2, lt2, accent1-4) plus the
# extended slots 9-12 (accent5, accent6, hlink, folHlink) that this
# color scheme exposes. RGB values are encoded the standard COM way
# (0x00BBGGRR -> decimal) since the RGB() helper isn't available here.
$scheme.Colors(1).RGB  = 0            # dk1      000000
$scheme.Colors(2).RGB  = 16777215     # lt1      FFFFFF
$scheme.Colors(3).RGB  = 6968388      # dk2      44546A
$scheme.Colors(4).RGB  = 15132391     # lt2      E7E6E6
$scheme.Colors(5).RGB  = 13998939     # accent1  5B9BD5
$scheme.Colors(6).RGB  = 3243501      # accent2  ED7D31
$scheme.Colors(7).RGB  = 10855845     # accent3  A5A5A5
$scheme.Colors(8).RGB  = 49407        # accent4  FFC000
$scheme.Colors(9).RGB  = 12874308     # accent5  4472C4
$scheme.Colors(10).RGB = 4697456      # accent6  70AD47
$scheme.Colors(11).RGB = 12673797     # hlink    0563C1
$scheme.Colors(12).RGB = 7491477      # folHlink 954F72
